$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Grants")

# Add new row 11 with label and sum formulas for B:E (sum of rows 7-9)
$ws.Range("A11").Value = "Other(Gift and adminstrative fees, Interest earnings granted, Cost recovery fees on investment goods)"
$ws.Range("B11").Formula = "=SUM(B7:B9)"
$ws.Range("C11").Formula = "=SUM(C7:C9)"
$ws.Range("D11").Formula = "=SUM(D7:D9)"
$ws.Range("E11").Formula = "=SUM(E7:E9)"

# Update the selection shown in the sheet view
$ws.Activate()
$ws.Range("A16").Select()

$wb.Save()
